$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full roster table after the edit (player, position, team), starting at row 2
$data = @(
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Vasilije Micic", "PG,SG", "Charlotte Hornets"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
